# Generate Report for Handback
# Updates timestamps (and one status value) on the handback-status report,
# simulating a later report-generation run.

$wb = $excel.ActiveWorkbook

# --- "Overview" sheet ---
# Column G = "Latest HO Xliff Generate Date"
# Rows 2 and 4 both previously showed 2016-08-19 16:14:34 -> 2016-08-19 16:15:25
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-08-19 16:15:25"
$wsOverview.Range("G4").Value = "2016-08-19 16:15:25"

# --- "zh-cn" sheet ---
# Column E = "Priority": row 2 and 4 go from "ht" to "mt"
# Column H = "Correspond Handoff Datetime": row 2 and 4 go from 16:14:30 to 16:15:21
# Column K = "Correspond Handback DateTime": row 2 and 4 go from 16:14:48 to 16:15:37
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "mt"
$wsZhCn.Range("E4").Value = "mt"
$wsZhCn.Range("H2").Value = "2016-08-19 16:15:21"
$wsZhCn.Range("H4").Value = "2016-08-19 16:15:21"
$wsZhCn.Range("K2").Value = "2016-08-19 16:15:37"
$wsZhCn.Range("K4").Value = "2016-08-19 16:15:37"

# --- "de-de" sheet ---
# Column H = "Correspond Handoff Datetime": row 2 and 4 go from 16:14:34 to 16:15:25
# (this reuses the same shared string as Overview!G2/G4 in the original workbook)
# Column K = "Correspond Handback DateTime": row 2 and 4 go from 16:14:55 to 16:15:44
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-08-19 16:15:25"
$wsDeDe.Range("H4").Value = "2016-08-19 16:15:25"
$wsDeDe.Range("K2").Value = "2016-08-19 16:15:44"
$wsDeDe.Range("K4").Value = "2016-08-19 16:15:44"
